$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text for column AC (29): drop the "Population: " prefix ---
$ws.Cells.Item(1, 29).Value = "Labor force participation (%)"

# --- Convert column AC data cells from text percentage strings (e.g. "48.90%")
#     to plain numeric values (e.g. 48.9) ---
$acData = @"
2,48.9
3,55.7
4,41.2
6,77.5
8,61.3
9,55.6
10,65.5
11,60.7
12,66.5
13,74.6
14,73.4
15,59
16,65.2
17,64.1
18,53.6
19,65.1
20,70.9
21,66.7
22,71.8
23,46.4
24,70.8
25,63.9
26,64.7
27,55.4
28,66.4
29,79.2
30,57
31,60.5
32,82.3
33,76.1
34,65.1
35,72
36,70.7
37,62.6
38,68
39,68.8
40,43.3
41,69.4
42,62.1
43,51.2
44,53.6
45,63.1
46,60.6
47,63.5
48,62.2
49,60.2
51,64.3
52,68
53,46.4
54,59.1
55,62
56,78.4
57,63.6
59,79.6
60,57.6
61,59.1
62,55.1
63,52.9
64,59.4
65,68.3
66,60.8
67,67.8
68,51.8
70,62.3
71,61.5
72,72
73,56.2
74,67.2
76,68.8
77,56.5
78,75
79,49.3
80,67.5
81,44.7
82,43
83,62.1
84,64
85,49.6
86,66
87,61.7
88,39.3
89,68.8
90,74.7
92,73.5
93,59.8
94,78.5
95,61.4
96,47
97,67.9
98,76.3
99,49.7
101,61.6
102,59.3
103,86.1
104,76.7
105,64.3
106,69.8
107,70.8
108,56.5
110,45.9
111,58.3
112,60.7
114,43.1
116,59.7
117,54.4
118,45.3
119,78.1
120,61.7
121,59.5
123,83.8
124,63.6
125,69.9
126,66.4
127,72
128,52.9
129,80.4
131,63.8
132,72.4
133,52.6
136,66.6
137,47.2
138,72.1
139,77.6
140,59.6
141,56.7
142,58.8
143,86.8
144,54.7
145,61.8
146,83.7
148,67.1
149,65.9
150,43.7
153,55.9
154,45.7
155,54.9
157,57.9
158,70.5
159,59.5
160,58.4
161,83.8
162,47.4
163,56
164,63
165,72.4
166,57.5
167,53.9
168,48.4
169,51.1
170,64.6
171,68.3
172,44.1
173,42
174,83.4
175,67.3
176,67.3
177,77.6
178,59.8
179,60
180,46.1
181,52.8
182,64.5
184,70.3
185,54.2
186,82.1
187,62.8
188,62
189,64
190,65.1
191,69.9
192,59.7
193,77.4
194,38
195,74.6
196,83.1
"@

$rows = $acData -split "`n"
foreach ($row in $rows) {
    $row = $row.Trim()
    if ($row -eq "") { continue }
    $parts = $row -split ","
    $rowNum = [int]$parts[0]
    $val = [double]$parts[1]
    $ws.Cells.Item($rowNum, 29).Value = $val
}
